# Update DM integration fixture hierarchies
# Replaces the ID (UUID) values on the CodeSchemes, Codes, and Extensions
# sheets with freshly generated UUIDs, matching the regenerated fixture data.

$wb = $excel.ActiveWorkbook

# --- CodeSchemes sheet ---
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A2").Value = "96ea0b5e-6a08-4a5b-b339-43cd5a0fbf81"

# --- Codes sheet ---
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A2").Value = "b81b6180-8b39-4edd-9fb9-3ef32622f0b0"
$wsCodes.Range("K2").Value = "d67ca944-01ab-4f4b-ba5f-0f2d6482c990"

$wsCodes.Range("A3").Value = "6ba70523-0b4f-45e0-bff0-d248c3824832"
$wsCodes.Range("K3").Value = "61d21a73-fc06-4697-8ebd-723509c1cc37"

$wsCodes.Range("A4").Value = "e19f54f2-bd3b-482a-8087-ea962ff9f548"
$wsCodes.Range("K4").Value = "fc765e19-55a7-4f7c-af48-028b10acda98"

$wsCodes.Range("A5").Value = "b761e37d-2700-4cff-aec8-728a87e0b31c"
$wsCodes.Range("K5").Value = "2ed27bd1-68a6-4724-bccb-46943652c845"

$wsCodes.Range("A6").Value = "13510de6-d515-4ed4-a2c7-6b65ff98edb6"
$wsCodes.Range("K6").Value = "96b8aa1f-84a9-4ffe-824d-4e1219e3f2a1"

$wsCodes.Range("A7").Value = "439a12ed-bfa3-4366-94f8-9b592b837929"
$wsCodes.Range("K7").Value = "b0603cfb-4121-405f-9ac8-540248a41235"

$wsCodes.Range("A8").Value = "de759e31-9cae-4e06-a8dd-0d8c8cab7faf"
$wsCodes.Range("K8").Value = "3713b47f-fdaa-4ae7-98ed-daeefd116a63"

# --- Extensions sheet ---
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A2").Value = "23ee2455-5c2a-45bb-991e-c924eaa67e72"

# Columns A (on all three sheets) use bestFit/customWidth so Excel
# recalculates the "best fit" width once the new (shorter/longer) UUID
# text is entered. Column K on the Codes sheet also narrows slightly as a
# side effect of the same autofit recompute.
$wsCodeSchemes.Columns.Item(1).ColumnWidth = 30.800000000000004
$wsCodes.Columns.Item(1).ColumnWidth = 35.2
$wsCodes.Columns.Item(11).ColumnWidth = 34.1
$wsExtensions.Columns.Item(1).ColumnWidth = 33.0
